$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume (E) columns to text format so that
# numeric-looking strings (e.g. "307.80", "1.004") are preserved verbatim
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '24.533.21'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '1.658.78'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '307.80'
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").Value = '0.3616'
$ws.Range("E7").Value = '  -3.15%  '
$ws.Range("D8").Value = '47.45'
$ws.Range("E8").Value = '  -3.00%  '
$ws.Range("D9").Value = '0.3247'
$ws.Range("E9").Value = '  -5.46%  '
$ws.Range("D10").Value = '1.122'
$ws.Range("E10").Value = '  -4.99%  '
$ws.Range("D11").Value = '0.06957'
$ws.Range("E11").Value = '  -6.62%  '
$ws.Range("D12").Value = '0.9999'
$ws.Range("D13").Value = '5.879'
$ws.Range("E13").Value = '  -5.71%  '
$ws.Range("D14").Value = '19.38'
$ws.Range("E14").Value = '  -7.13%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.550'
$ws.Range("E15").Value = '  -5.42%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.649.68'
$ws.Range("E16").Value = '  -3.32%  '
$ws.Range("D17").Value = '0.00001044'
$ws.Range("E17").Value = '  -6.65%  '
$ws.Range("D18").Value = '0.06544'
$ws.Range("E18").Value = '  -2.31%  '
$ws.Range("D19").Value = '0.9992'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '76.37'
$ws.Range("E20").Value = '  -8.64%  '
$ws.Range("D21").Value = '5.912'
$ws.Range("E21").Value = '  -6.50%  '
$ws.Range("D22").Value = '15.65'
$ws.Range("E22").Value = '  -8.37%  '
$ws.Range("D23").Value = '12.62'
$ws.Range("E23").Value = '  -3.70%  '
$ws.Range("D24").Value = '24.491.12'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").Value = '2.465'
$ws.Range("E25").Value = '  +2.26%  '
$ws.Range("D26").Value = '2.297'
$ws.Range("E26").Value = '  -16.71%  '
$ws.Range("D27").Value = '146.83'
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").Value = '18.42'
$ws.Range("E28").Value = '  -8.51%  '
$ws.Range("D29").Value = '1.834.96'
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("D30").Value = '1.193'
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("D31").Value = '123.75'
$ws.Range("E31").Value = '  -5.69%  '
$ws.Range("D32").Value = '4.067'
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("D33").Value = '5.620'
$ws.Range("E33").Value = '  -16.61%  '
$ws.Range("D34").Value = '1.699'
$ws.Range("E34").Value = '  -4.29%  '
$ws.Range("D35").Value = '0.08347'
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("D36").Value = '12.34'
$ws.Range("E36").Value = '  -9.61%  '
$ws.Range("D37").Value = '5.192'
$ws.Range("E37").Value = '  -5.81%  '
$ws.Range("D38").Value = '0.06041'
$ws.Range("E38").Value = '  -7.13%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02188'
$ws.Range("E39").Value = '  -7.88%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.212'
$ws.Range("E40").Value = '  -7.92%  '
$ws.Range("D41").Value = '0.2052'
$ws.Range("E41").Value = '  -7.41%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.201'
$ws.Range("E42").Value = '  -5.75%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").Value = '0.5889'
$ws.Range("E44").Value = '  -7.76%  '
$ws.Range("D45").Value = '3.738'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").Value = '12.69'
$ws.Range("E46").Value = '  -8.65%  '
$ws.Range("D47").Value = '0.5581'
$ws.Range("E47").Value = '  -8.14%  '
$ws.Range("D48").Value = '122.37'
$ws.Range("E48").Value = '  -5.29%  '
$ws.Range("D49").Value = '1.937'
$ws.Range("E49").Value = '  -8.25%  '
$ws.Range("D50").Value = '0.06909'
$ws.Range("E50").Value = '  -4.91%  '
$ws.Range("D51").Value = '74.03'
$ws.Range("E51").Value = '  -6.35%  '
